$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.46025589921483
$ws.Range("C2").Value = -0.3606512508042645
$ws.Range("D2").Value = -17.46025589921483
$ws.Range("E2").Value = -17.46025589921483
$ws.Range("F2").Value = -17.46025589921483
$ws.Range("G2").Value = -17.46025589921483
$ws.Range("H2").Value = -17.46025589921483
$ws.Range("I2").Value = -17.46025589921483
$ws.Range("J2").Value = -17.46025589921483
$ws.Range("K2").Value = -17.46025589921483
$ws.Range("B3").Value = -17.46025589921483
$ws.Range("C3").Value = -17.46025589921483
$ws.Range("D3").Value = -17.46025589921483
$ws.Range("E3").Value = -17.46025589921483
$ws.Range("F3").Value = -17.46025589921483
$ws.Range("G3").Value = -17.46025589921483
$ws.Range("H3").Value = -17.46025589921483
$ws.Range("I3").Value = 0.4067633656820623
$ws.Range("J3").Value = -17.46025589921483
$ws.Range("K3").Value = -17.46025589921483
$ws.Range("B4").Value = -17.46025589921483
$ws.Range("C4").Value = -0.5818097751447958
$ws.Range("D4").Value = 0.4919450095486629
$ws.Range("E4").Value = -17.46025589921483
$ws.Range("F4").Value = 4.068028887779802
$ws.Range("G4").Value = -17.46025589921483
$ws.Range("H4").Value = 1.588903994568738
$ws.Range("I4").Value = -17.46025589921483
$ws.Range("J4").Value = 2.824004538161883
$ws.Range("K4").Value = -17.46025589921483
$ws.Range("B5").Value = -17.46025589921483
$ws.Range("C5").Value = 0.135183714483207
$ws.Range("D5").Value = -17.46025589921483
$ws.Range("E5").Value = -17.46025589921483
$ws.Range("F5").Value = -17.46025589921483
$ws.Range("G5").Value = 3.488579695299702
$ws.Range("H5").Value = -17.46025589921483
$ws.Range("I5").Value = -17.46025589921483
$ws.Range("J5").Value = -17.46025589921483
$ws.Range("K5").Value = -17.46025589921483
$ws.Range("B6").Value = -17.46025589921483
$ws.Range("C6").Value = -17.46025589921483
$ws.Range("D6").Value = -17.46025589921483
$ws.Range("E6").Value = -17.46025589921483
$ws.Range("F6").Value = -17.46025589921483
$ws.Range("G6").Value = -17.46025589921483
$ws.Range("H6").Value = -17.46025589921483
$ws.Range("I6").Value = -17.46025589921483
$ws.Range("J6").Value = -17.46025589921483
$ws.Range("K6").Value = -17.46025589921483
$ws.Range("B7").Value = 3.027777401418363
$ws.Range("C7").Value = -17.46025589921483
$ws.Range("D7").Value = -17.46025589921483
$ws.Range("E7").Value = -17.46025589921483
$ws.Range("F7").Value = -17.46025589921483
$ws.Range("G7").Value = -17.46025589921483
$ws.Range("H7").Value = -17.46025589921483
$ws.Range("I7").Value = -17.46025589921483
$ws.Range("J7").Value = -17.46025589921483
$ws.Range("K7").Value = -17.46025589921483
$ws.Range("B8").Value = -17.46025589921483
$ws.Range("C8").Value = -17.46025589921483
$ws.Range("D8").Value = -17.46025589921483
$ws.Range("E8").Value = -17.46025589921483
$ws.Range("F8").Value = -17.46025589921483
$ws.Range("G8").Value = -17.46025589921483
$ws.Range("H8").Value = -17.46025589921483
$ws.Range("I8").Value = -17.46025589921483
$ws.Range("J8").Value = -17.46025589921483
$ws.Range("K8").Value = -17.46025589921483
$ws.Range("B9").Value = 3.566130430835261
$ws.Range("C9").Value = -17.46025589921483
$ws.Range("D9").Value = -17.46025589921483
$ws.Range("E9").Value = -17.46025589921483
$ws.Range("F9").Value = -17.46025589921483
$ws.Range("G9").Value = -17.46025589921483
$ws.Range("H9").Value = -17.46025589921483
$ws.Range("I9").Value = -17.46025589921483
$ws.Range("J9").Value = -17.46025589921483
$ws.Range("K9").Value = -17.46025589921483
$ws.Range("B10").Value = -17.46025589921483
$ws.Range("C10").Value = -17.46025589921483
$ws.Range("D10").Value = -17.46025589921483
$ws.Range("E10").Value = -17.46025589921483
$ws.Range("F10").Value = -17.46025589921483
$ws.Range("G10").Value = -17.46025589921483
$ws.Range("H10").Value = -17.46025589921483
$ws.Range("I10").Value = 0.2383470396458971
$ws.Range("J10").Value = -17.46025589921483
$ws.Range("K10").Value = 2.000745848044908
$ws.Range("B11").Value = -17.46025589921483
$ws.Range("C11").Value = -17.46025589921483
$ws.Range("D11").Value = -17.46025589921483
$ws.Range("E11").Value = -17.46025589921483
$ws.Range("F11").Value = -17.46025589921483
$ws.Range("G11").Value = 1.349703730533285
$ws.Range("H11").Value = -17.46025589921483
$ws.Range("I11").Value = -17.46025589921483
$ws.Range("J11").Value = -17.46025589921483
$ws.Range("K11").Value = 1.167951099192089
$ws.Range("B12").Value = -17.46025589921483
$ws.Range("C12").Value = -17.46025589921483
$ws.Range("D12").Value = -17.46025589921483
$ws.Range("E12").Value = -17.46025589921483
$ws.Range("F12").Value = -17.46025589921483
$ws.Range("G12").Value = -17.46025589921483
$ws.Range("H12").Value = -17.46025589921483
$ws.Range("I12").Value = -17.46025589921483
$ws.Range("J12").Value = -17.46025589921483
$ws.Range("K12").Value = -17.46025589921483
$ws.Range("B13").Value = -17.46025589921483
$ws.Range("C13").Value = -17.46025589921483
$ws.Range("D13").Value = -17.46025589921483
$ws.Range("E13").Value = 4.321920494450938
$ws.Range("F13").Value = -17.46025589921483
$ws.Range("G13").Value = -17.46025589921483
$ws.Range("H13").Value = -17.46025589921483
$ws.Range("I13").Value = -17.46025589921483
$ws.Range("J13").Value = 0.9519541882167023
$ws.Range("K13").Value = 2.740323068496634
$ws.Range("B14").Value = -17.46025589921483
$ws.Range("C14").Value = -17.46025589921483
$ws.Range("D14").Value = -0.24255373908138
$ws.Range("E14").Value = -17.46025589921483
$ws.Range("F14").Value = -17.46025589921483
$ws.Range("G14").Value = -17.46025589921483
$ws.Range("H14").Value = -17.46025589921483
$ws.Range("I14").Value = -17.46025589921483
$ws.Range("J14").Value = -17.46025589921483
$ws.Range("K14").Value = 1.928604366118605
$ws.Range("B15").Value = -17.46025589921483
$ws.Range("C15").Value = -17.46025589921483
$ws.Range("D15").Value = -0.5470275358154805
$ws.Range("E15").Value = -17.46025589921483
$ws.Range("F15").Value = -17.46025589921483
$ws.Range("G15").Value = -17.46025589921483
$ws.Range("H15").Value = -17.46025589921483
$ws.Range("I15").Value = -17.46025589921483
$ws.Range("J15").Value = -17.46025589921483
$ws.Range("K15").Value = -17.46025589921483
$ws.Range("B16").Value = -17.46025589921483
$ws.Range("C16").Value = -17.46025589921483
$ws.Range("D16").Value = -17.46025589921483
$ws.Range("E16").Value = -17.46025589921483
$ws.Range("F16").Value = -17.46025589921483
$ws.Range("G16").Value = -17.46025589921483
$ws.Range("H16").Value = -17.46025589921483
$ws.Range("I16").Value = -17.46025589921483
$ws.Range("J16").Value = 1.932591853332147
$ws.Range("K16").Value = -17.46025589921483
$ws.Range("B17").Value = -17.46025589921483
$ws.Range("C17").Value = 1.140125283857859
$ws.Range("D17").Value = -0.2975400477844491
$ws.Range("E17").Value = -17.46025589921483
$ws.Range("F17").Value = -17.46025589921483
$ws.Range("G17").Value = -17.46025589921483
$ws.Range("H17").Value = 1.998745385925073
$ws.Range("I17").Value = 0.5221194169824812
$ws.Range("J17").Value = 1.89389230661108
$ws.Range("K17").Value = -17.46025589921483
$ws.Range("B18").Value = -17.46025589921483
$ws.Range("C18").Value = -17.46025589921483
$ws.Range("D18").Value = -17.46025589921483
$ws.Range("E18").Value = -17.46025589921483
$ws.Range("F18").Value = -17.46025589921483
$ws.Range("G18").Value = -17.46025589921483
$ws.Range("H18").Value = 2.255131367205801
$ws.Range("I18").Value = -0.2787833120097445
$ws.Range("J18").Value = 1.78676528879226
$ws.Range("K18").Value = -17.46025589921483
$ws.Range("B19").Value = -17.46025589921483
$ws.Range("C19").Value = -17.46025589921483
$ws.Range("D19").Value = 3.25138839092096
$ws.Range("E19").Value = -17.46025589921483
$ws.Range("F19").Value = -17.46025589921483
$ws.Range("G19").Value = -17.46025589921483
$ws.Range("H19").Value = 1.976839832274682
$ws.Range("I19").Value = 1.205832844587907
$ws.Range("J19").Value = -17.46025589921483
$ws.Range("K19").Value = -17.46025589921483
$ws.Range("B20").Value = -17.46025589921483
$ws.Range("C20").Value = 3.216679957281436
$ws.Range("D20").Value = 2.75005991073994
$ws.Range("E20").Value = -17.46025589921483
$ws.Range("F20").Value = 1.690356629295031
$ws.Range("G20").Value = -17.46025589921483
$ws.Range("H20").Value = 1.07478021347457
$ws.Range("I20").Value = 3.692373694166331
$ws.Range("J20").Value = -17.46025589921483
$ws.Range("K20").Value = 1.705693812173613
$ws.Range("B21").Value = -17.46025589921483
$ws.Range("C21").Value = 2.573931836424259
$ws.Range("D21").Value = -17.46025589921483
$ws.Range("E21").Value = -17.46025589921483
$ws.Range("F21").Value = -17.46025589921483
$ws.Range("G21").Value = 2.638492016349776
$ws.Range("H21").Value = 1.123476974964175
$ws.Range("I21").Value = -17.46025589921483
$ws.Range("J21").Value = -17.46025589921483
$ws.Range("K21").Value = -17.46025589921483
